$wb = $excel.ActiveWorkbook

# =================================================================
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet
# =================================================================
$totalBefore = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2021-Q4")

$wb.Worksheets.Add($totalBefore) | Out-Null

# NOTE: re-fetch both sheets by name after Add() - the object returned by
# Add()/passed as "Before" can end up referring to the freshly inserted
# sheet once it is renamed, so always look sheets up fresh by name.
$ws = $wb.Worksheets.Item(5)
$ws.Name = "2022-Q1"

$total = $wb.Worksheets.Item("总计")

# Copy cell formatting (font/border/alignment) from the template sheet so
# the new sheet matches the look of the other quarterly sheets. The new
# sheet needs exactly 7 rows (1 header + 6 data rows). Column A on the
# header row is intentionally left untouched (the template has no A1 cell
# either).
$template.Range("B1:H7").Copy()
$ws.Range("B1:H7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$template.Range("A2:A7").Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# a never-touched cell used as a "blank format" source, to strip away the
# temporary text NumberFormat we apply below once values are in place.
$blank = $ws.Range("Z1")

# ---------------- Header row -----------------
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# ---------------- Data rows ------------------
# columns: seq, code, name, scale, stockPosition, positionRatio, marketValue, rank
$rows = @(
    @(0, "008704", "广发高股息优享混合A", "3.52", "94.09", "5.55", "0.1954", 7),
    @(1, "008705", "广发高股息优享混合C", "0.82", "94.09", "5.55", "0.0455", 7),
    @(2, "010756", "兴华永兴混合A",       "0.35", "94.57", "3.70", "0.0130", 9),
    @(3, "010999", "兴华瑞丰混合A",       "0.06", "29.21", "3.39", "0.0020", 5),
    @(4, "011000", "兴华瑞丰混合C",       "0.05", "29.21", "3.39", "0.0017", 5),
    @(5, "010757", "兴华永兴混合C",       "0.01", "94.57", "3.70", "0.0004", 9)
)

$r = 2
foreach ($item in $rows) {
    $ws.Cells.Item($r, 1).Value = $item[0]

    # B, C, D, E, F, G must stay as *text* (B/C look like they could be
    # coerced to numbers because of the leading zeros, D..G look numeric
    # too) - force text format before writing the value, then strip the
    # number-format residue by pasting a blank cell's formatting over it.
    $textCells = $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 7))
    $textCells.NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    $ws.Cells.Item($r, 4).Value = $item[3]
    $ws.Cells.Item($r, 5).Value = $item[4]
    $ws.Cells.Item($r, 6).Value = $item[5]
    $ws.Cells.Item($r, 7).Value = $item[6]
    $blank.Copy()
    $textCells.PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Cells.Item($r, 8).Value = $item[7]
    $r = $r + 1
}

# =================================================================
# 2. Update the "总计" sheet: prepend a new row for 2022-Q1
# =================================================================
$total.Rows.Item(2).Insert(-4121)
$total.Rows.Item(2).ClearFormats()

# restore formatting of the new A2 cell to match the other sequence cells
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 6
$total.Cells.Item(2, 4).Value = 0.26

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(6, 1).Value = 4
